$d = $word.ActiveDocument

# First paragraph holds the document title: "Jenkins with Terraform Documentation"
$p = $d.Paragraphs(1)

# Indent the title 1 inch (1440 twips = 72 points) from the left margin.
$p.Format.LeftIndent = 72

# Make the title bold and bump it up to 18pt (36 half-points), for both the
# normal and the complex-script (bidi) font properties, so both
# <w:b/><w:bCs/> and <w:sz/><w:szCs/> are written out.
$p.Range.Font.Bold = $true
$p.Range.Font.BoldBi = $true
$p.Range.Font.Size = 18
$p.Range.Font.SizeBi = 18
